# Rename the PDF file names listed in column C (rows 1-14) to the new
# normalized naming scheme (no accents/underscores/dots, numbered, fixed
# ordering between "AutoAdmiteDemanda" and "ActaDeReparto").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newNames = @(
    "01Caratula.pdf",
    "02AutoAdmiteDemanda.pdf",
    "03ActaDeReparto.pdf",
    "04Memorial.pdf",
    "05Constancia20210321.pdf",
    "06NotificacionDemandado.pdf",
    "07MemorialNoAceptaDesignacion.pdf",
    "08AutoNombraCurador.pdf",
    "09AcuseRecibido.pdf",
    "10AceptaDesignacion.pdf",
    "11NotificacionCurador.pdf",
    "12ConstestacionCuradorAdLitem.pdf",
    "13MemorialSolicitudCeleridad.pdf",
    "14AutoOrdenaSeguirAdelanteEjecucion.pdf"
)

for ($i = 0; $i -lt $newNames.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $newNames[$i]
}
